$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Title / header text updates ---
$ws.Range("A8").Value = "Volume 33   Number  2"
$ws.Range("C9").Value = "Report Covering the Week  1/5/2026  Through  1/11/2026"

# --- Template cells used to correctly write literal-text "0" / "***.*" values ---
$zeroTemplate = $ws.Range("C31")   # contains literal text "0" (style General)
$naTemplate   = $ws.Range("E31")   # contains literal text "***.*" (style General)

# Row 14
$ws.Range("D14").Value = 1
$ws.Range("G14").Value = 3
$ws.Range("H14").Value = -66.666666666666
$ws.Range("J14").Value = 3
$ws.Range("L14").Value = -100
$ws.Range("L14").NumberFormat = "#,##0.0;""-""#,##0.0"

# Row 15
$zeroTemplate.Copy() | Out-Null
$ws.Range("C15").PasteSpecial(-4163) | Out-Null   # values (brings literal text)
$zeroTemplate.Copy() | Out-Null
$ws.Range("C15").PasteSpecial(-4122) | Out-Null   # formats (restore General style)
$ws.Range("D15").Value = 1
$ws.Range("D15").NumberFormat = "#,##0"
$ws.Range("E15").Value = -100
$ws.Range("E15").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("J15").Value = 1
$ws.Range("J15").NumberFormat = "#,##0"
$ws.Range("K15").Value = 0
$ws.Range("K15").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("L15").Value = 0
$ws.Range("L15").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("N15").Value = 0
$ws.Range("N15").NumberFormat = "#,##0.0;""-""#,##0.0"

# Row 16
$ws.Range("C16").Value = 5
$ws.Range("D16").Value = 4
$ws.Range("E16").Value = 25
$ws.Range("F16").Value = 29
$ws.Range("G16").Value = 27
$ws.Range("H16").Value = 7.407407407407
$ws.Range("I16").Value = 10
$ws.Range("J16").Value = 10
$ws.Range("K16").Value = 0
$ws.Range("L16").Value = -56.521739130434
$ws.Range("M16").Value = 11.111111111111
$ws.Range("N16").Value = -87.80487804878

# Row 17
$ws.Range("C17").Value = 14
$ws.Range("D17").Value = 16
$ws.Range("E17").Value = -12.5
$ws.Range("F17").Value = 57
$ws.Range("G17").Value = 82
$ws.Range("H17").Value = -30.487804878048
$ws.Range("I17").Value = 20
$ws.Range("J17").Value = 28
$ws.Range("K17").Value = -28.571428571428
$ws.Range("L17").Value = -37.5
$ws.Range("M17").Value = 66.666666666666
$ws.Range("N17").Value = -42.857142857142

# Row 18
$ws.Range("C18").Value = 2
$ws.Range("D18").Value = 3
$ws.Range("E18").Value = -33.333333333333
$ws.Range("F18").Value = 13
$ws.Range("G18").Value = 20
$ws.Range("H18").Value = -35
$ws.Range("I18").Value = 6
$ws.Range("J18").Value = 3
$ws.Range("J18").NumberFormat = "#,##0"
$ws.Range("K18").Value = 100
$ws.Range("K18").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("L18").Value = -33.333333333333
$ws.Range("M18").Value = -33.333333333333
$ws.Range("N18").Value = -92.105263157894

# Row 19
$ws.Range("C19").Value = 14
$ws.Range("D19").Value = 14
$ws.Range("E19").Value = 0
$ws.Range("F19").Value = 50
$ws.Range("G19").Value = 50
$ws.Range("H19").Value = 0
$ws.Range("I19").Value = 18
$ws.Range("J19").Value = 22
$ws.Range("K19").Value = -18.181818181818
$ws.Range("L19").Value = -25
$ws.Range("M19").Value = 200
$ws.Range("N19").Value = 12.5

# Row 20
$ws.Range("C20").Value = 5
$ws.Range("E20").Value = 66.666666666666
$ws.Range("F20").Value = 25
$ws.Range("G20").Value = 14
$ws.Range("H20").Value = 78.571428571428
$ws.Range("I20").Value = 8
$ws.Range("J20").Value = 5
$ws.Range("K20").Value = 60
$ws.Range("L20").Value = 0
$ws.Range("M20").Value = 166.666666666667
$ws.Range("N20").Value = -83.673469387755

# Row 21
$ws.Range("D21").Value = 42
$ws.Range("E21").Value = -4.761904761904
$ws.Range("F21").Value = 178
$ws.Range("G21").Value = 198
$ws.Range("H21").Value = -10.10101010101
$ws.Range("I21").Value = 63
$ws.Range("J21").Value = 72
$ws.Range("K21").Value = -12.5
$ws.Range("L21").Value = -35.714285714285
$ws.Range("M21").Value = 61.538461538461
$ws.Range("N21").Value = -75.954198473282

# Row 22
$zeroTemplate.Copy() | Out-Null
$ws.Range("C22").PasteSpecial(-4163) | Out-Null   # values (brings literal text)
$zeroTemplate.Copy() | Out-Null
$ws.Range("C22").PasteSpecial(-4122) | Out-Null   # formats (restore General style)
$ws.Range("D22").Value = 1
$ws.Range("D22").NumberFormat = "#,##0"
$ws.Range("E22").Value = -100
$ws.Range("E22").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("F22").Value = 2
$ws.Range("G22").Value = 2
$ws.Range("H22").Value = 0
$ws.Range("J22").Value = 1
$ws.Range("J22").NumberFormat = "#,##0"
$ws.Range("K22").Value = 0
$ws.Range("K22").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("L22").Value = -75
$ws.Range("M22").Value = 0
$ws.Range("M22").NumberFormat = "#,##0.0;""-""#,##0.0"

# Row 23
$zeroTemplate.Copy() | Out-Null
$ws.Range("C23").PasteSpecial(-4163) | Out-Null   # values (brings literal text)
$zeroTemplate.Copy() | Out-Null
$ws.Range("C23").PasteSpecial(-4122) | Out-Null   # formats (restore General style)
$zeroTemplate.Copy() | Out-Null
$ws.Range("D23").PasteSpecial(-4163) | Out-Null   # values (brings literal text)
$zeroTemplate.Copy() | Out-Null
$ws.Range("D23").PasteSpecial(-4122) | Out-Null   # formats (restore General style)
$naTemplate.Copy() | Out-Null
$ws.Range("E23").PasteSpecial(-4163) | Out-Null   # values (brings literal text)
$naTemplate.Copy() | Out-Null
$ws.Range("E23").PasteSpecial(-4122) | Out-Null   # formats (restore General style)
$ws.Range("F23").Value = 5
$ws.Range("H23").Value = 0
$ws.Range("L23").Value = 0
$ws.Range("L23").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("M23").Value = 100
$ws.Range("M23").NumberFormat = "#,##0.0;""-""#,##0.0"

# Row 24
$ws.Range("C24").Value = 22
$ws.Range("D24").Value = 41
$ws.Range("E24").Value = -46.341463414634
$ws.Range("F24").Value = 102
$ws.Range("G24").Value = 155
$ws.Range("H24").Value = -34.193548387096
$ws.Range("I24").Value = 32
$ws.Range("J24").Value = 54
$ws.Range("K24").Value = -40.74074074074
$ws.Range("L24").Value = -40.74074074074
$ws.Range("M24").Value = 10.344827586206

# Row 25
$ws.Range("C25").Value = 7
$ws.Range("D25").Value = 16
$ws.Range("E25").Value = -56.25
$ws.Range("F25").Value = 27
$ws.Range("G25").Value = 47
$ws.Range("H25").Value = -42.553191489361
$ws.Range("I25").Value = 11
$ws.Range("J25").Value = 22
$ws.Range("K25").Value = -50
$ws.Range("L25").Value = -67.647058823529

# Row 26
$ws.Range("C26").Value = 27
$ws.Range("D26").Value = 23
$ws.Range("E26").Value = 17.391304347826
$ws.Range("F26").Value = 99
$ws.Range("G26").Value = 102
$ws.Range("H26").Value = -2.941176470588
$ws.Range("I26").Value = 38
$ws.Range("J26").Value = 30
$ws.Range("K26").Value = 26.666666666666
$ws.Range("L26").Value = 31.03448275862
$ws.Range("M26").Value = 111.111111111111

# Row 27
$zeroTemplate.Copy() | Out-Null
$ws.Range("C27").PasteSpecial(-4163) | Out-Null   # values (brings literal text)
$zeroTemplate.Copy() | Out-Null
$ws.Range("C27").PasteSpecial(-4122) | Out-Null   # formats (restore General style)
$ws.Range("D27").Value = 3
$ws.Range("D27").NumberFormat = "#,##0"
$ws.Range("E27").Value = -100
$ws.Range("E27").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("G27").Value = 5
$ws.Range("H27").Value = -40
$ws.Range("J27").Value = 3
$ws.Range("J27").NumberFormat = "#,##0"
$ws.Range("K27").Value = -66.666666666666
$ws.Range("K27").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("L27").Value = 0
$ws.Range("L27").NumberFormat = "#,##0.0;""-""#,##0.0"

# Row 28
$ws.Range("C28").Value = 1
$ws.Range("D28").Value = 3
$ws.Range("E28").Value = -66.666666666666
$ws.Range("F28").Value = 7
$ws.Range("G28").Value = 11
$ws.Range("H28").Value = -36.363636363636
$ws.Range("I28").Value = 3
$ws.Range("J28").Value = 5
$ws.Range("K28").Value = -40
$ws.Range("L28").Value = 50
$ws.Range("L28").NumberFormat = "#,##0.0;""-""#,##0.0"

# Row 29
$zeroTemplate.Copy() | Out-Null
$ws.Range("C29").PasteSpecial(-4163) | Out-Null   # values (brings literal text)
$zeroTemplate.Copy() | Out-Null
$ws.Range("C29").PasteSpecial(-4122) | Out-Null   # formats (restore General style)
$zeroTemplate.Copy() | Out-Null
$ws.Range("D29").PasteSpecial(-4163) | Out-Null   # values (brings literal text)
$zeroTemplate.Copy() | Out-Null
$ws.Range("D29").PasteSpecial(-4122) | Out-Null   # formats (restore General style)
$naTemplate.Copy() | Out-Null
$ws.Range("E29").PasteSpecial(-4163) | Out-Null   # values (brings literal text)
$naTemplate.Copy() | Out-Null
$ws.Range("E29").PasteSpecial(-4122) | Out-Null   # formats (restore General style)
$ws.Range("M29").Value = 200
$ws.Range("M29").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("N29").Value = 50

# Row 30
$zeroTemplate.Copy() | Out-Null
$ws.Range("C30").PasteSpecial(-4163) | Out-Null   # values (brings literal text)
$zeroTemplate.Copy() | Out-Null
$ws.Range("C30").PasteSpecial(-4122) | Out-Null   # formats (restore General style)
$zeroTemplate.Copy() | Out-Null
$ws.Range("D30").PasteSpecial(-4163) | Out-Null   # values (brings literal text)
$zeroTemplate.Copy() | Out-Null
$ws.Range("D30").PasteSpecial(-4122) | Out-Null   # formats (restore General style)
$naTemplate.Copy() | Out-Null
$ws.Range("E30").PasteSpecial(-4163) | Out-Null   # values (brings literal text)
$naTemplate.Copy() | Out-Null
$ws.Range("E30").PasteSpecial(-4122) | Out-Null   # formats (restore General style)
$ws.Range("M30").Value = 100
$ws.Range("M30").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("N30").Value = 0

# Row 31
$zeroTemplate.Copy() | Out-Null
$ws.Range("F31").PasteSpecial(-4163) | Out-Null   # values (brings literal text)
$zeroTemplate.Copy() | Out-Null
$ws.Range("F31").PasteSpecial(-4122) | Out-Null   # formats (restore General style)

# --- Historical perspective totals (rows 45-46) ---
$ws.Range("J45").Value = 318
$ws.Range("K45").Value = -32.627118644067
$ws.Range("L45").Value = -38.728323699422
$ws.Range("M45").Value = -75.136825645035
$ws.Range("N45").Value = -81.766055045871

$ws.Range("J46").Value = 2834
$ws.Range("K46").Value = -1.971636112071
$ws.Range("L46").Value = -22.568306010929
$ws.Range("M46").Value = -63.624695161083
$ws.Range("N46").Value = -67.850255246738
